$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 337, shifting existing rows 337:364 down to 338:365
$ws.Rows.Item(337).Insert()

# Populate the newly inserted row 337 with the new weekly data point
$ws.Cells.Item(337, 1).Value = 8
$ws.Cells.Item(337, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(337, 3).Value = "Coquimbo"
$ws.Cells.Item(337, 4).Value = 45013
$ws.Cells.Item(337, 5).Value = 4
$ws.Cells.Item(337, 6).Value = 100112031
$ws.Cells.Item(337, 7).Value = "Poroto verde"
$ws.Cells.Item(337, 8).Value = "Magnum"
$ws.Cells.Item(337, 9).Value = "Primera"
$ws.Cells.Item(337, 10).Value = 400
$ws.Cells.Item(337, 11).Value = 22000
$ws.Cells.Item(337, 12).Value = 23000
$ws.Cells.Item(337, 13).Value = 22500
$ws.Cells.Item(337, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(337, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(337, 16).Value = 900
$ws.Cells.Item(337, 17).Value = 25
$ws.Cells.Item(337, 18).Value = "Hortaliza"
